$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. TABLE_NAMES sheet: append the new table name "T_EMP" as the next row.
# ---------------------------------------------------------------------------
$wsTableNames = $wb.Worksheets.Item("TABLE_NAMES")
$wsTableNames.Range("A5").Value = "T_EMP"
$wsTableNames.Range("A5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. STUDENT sheet used to be the active/selected tab - move the selection to
#    the header row; the new sheet will become the active tab instead.
# ---------------------------------------------------------------------------
$wsStudent = $wb.Worksheets.Item("STUDENT")
$wsStudent.Range("A1:C1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Add the new T_EMP worksheet after the last existing sheet (STUDENT).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTEmp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsTEmp.Name = "T_EMP"

# Header row.
$wsTEmp.Range("A1").Value = "COULMN_NAME"
$wsTEmp.Range("B1").Value = "DATA_TYPE"
$wsTEmp.Range("C1").Value = "DATA_LENGTH"

# Column metadata rows for the new T_EMP table.
$rows = @(
    @("T_EMP_ID",              "int",     "int"),
    @("EMP_ID",                "int",     "int"),
    @("EMP_NAME",               "varchar", "varchar(50)"),
    @("DEPT_NAME",              "varchar", "varchar(20)"),
    @("MANAGER_NAME",           "varchar", "varchar(50)"),
    @("JOB_TYPE",                "varchar", "varchar(20)"),
    @("LOCATION",                 "varchar", "varchar(30)"),
    @("HIRE_DATE",               "date",    "date"),
    @("TOTAL_EXP_IN_COMPANY",   "int",     "int"),
    @("SALARY",                  "int",     "int"),
    @("COMMISSION",              "int",     "int"),
    @("TOTAL_SALARY",           "int",     "int"),
    @("SALARY_GRADE",           "varchar", "varchar(10)")
)

$r = 2
foreach ($row in $rows) {
    $wsTEmp.Cells.Item($r, 1).Value = $row[0]
    $wsTEmp.Cells.Item($r, 2).Value = $row[1]
    $wsTEmp.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Match the formatting used on the other table sheets (wrap text, vertically
# centred) for the data rows beneath the header - copy it from an already
# formatted sheet so we reuse the existing cell style instead of minting a
# new one.
$wsStudent.Range("A2:C2").Copy() | Out-Null
$wsTEmp.Range("A2:C14").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

# Size the columns to fit their content, like the other sheets in the book.
$wsTEmp.Columns.Item(1).AutoFit() | Out-Null
$wsTEmp.Columns.Item(2).AutoFit() | Out-Null
$wsTEmp.Columns.Item(3).AutoFit() | Out-Null

# Leave the cursor where the author left off and make T_EMP the active tab.
$wsTEmp.Range("E10").Select() | Out-Null
$wsTEmp.Activate() | Out-Null
